$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '25.909.61'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.75%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.631.27'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.09%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.15%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.49'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.07%  '
$ws.Range("E6").Value = '  +0.37%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.20%  '
$ws.Range("E8").Value = '  +0.21%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0632'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.08%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.64'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.91%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0790'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.05%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.855.93'
$ws.Range("D12").Style = "Normal"
$ws.Range("E13").Value = '  -0.62%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.623.94'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.43%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.545'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.54%  '
$ws.Range("E16").Value = '  -0.46%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '62.78'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.03%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '25.885.74'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.63%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.00'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.14%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '193.27'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.89%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.95'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.30%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.27'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.38%  '
$ws.Range("E24").Value = '  -1.50%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.999'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.28%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '142.22'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.36%  '
$ws.Range("E27").Value = '  +2.49%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.86'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.46%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.46'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.07%  '
$ws.Range("E30").Value = '  -0.07%  '
$ws.Range("E31").Value = '  +1.82%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.32'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.23%  '
$ws.Range("E33").Value = '  -0.24%  '
$ws.Range("E34").Value = '  -0.22%  '
$ws.Range("E35").Value = '  +1.40%  '
$ws.Range("E36").Value = '  -0.27%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.134.16'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.75%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.551'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.89%  '
$ws.Range("E39").Value = '  -1.80%  '
$ws.Range("E40").Value = '  +0.44%  '
$ws.Range("E41").Value = '  -0.14%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.805'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.06%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '99.20'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.47%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.44'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.68%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.765.35'
$ws.Range("D45").Style = "Normal"
$ws.Range("E46").Value = '  -0.68%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '56.14'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.76%  '
$ws.Range("E48").Value = '  +3.53%  '
$ws.Range("E49").Value = '  +1.47%  '
$ws.Range("E50").Value = '  -1.00%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.60'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.32%  '
